$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 18.5
$ws.Range("I8").Value = 18.5
$ws.Range("K8").Value = 55.5
$ws.Range("M8").Value = 83.5
$ws.Range("H19").Value = 1121.7273
$ws.Range("I19").Value = 1146
$ws.Range("J19").Value = 1092.6
$ws.Range("K19").Value = 1146
$ws.Range("L19").Value = 1092.6
$ws.Range("M19").Value = -971
$ws.Range("N19").Value = -1442.6
$ws.Range("H33").Value = 8728.416999999999
$ws.Range("I33").Value = 11516.777
$ws.Range("J33").Value = 363.33334
$ws.Range("K33").Value = 11516.777
$ws.Range("L33").Value = 363.33334
$ws.Range("M33").Value = -11287.777
$ws.Range("N33").Value = -821.33334
$ws.Range("H86").Value = 5242.4546
$ws.Range("I86").Value = 4147.125
$ws.Range("K86").Value = 4147.125
$ws.Range("M86").Value = -3024.125
$ws.Range("H89").Value = 5242.4546
$ws.Range("I89").Value = 4147.125
$ws.Range("K89").Value = 20735.625
$ws.Range("M89").Value = -15119.625
$ws.Range("H135").Value = 1418.7
$ws.Range("I135").Value = 1414.1154
$ws.Range("J135").Value = 1448.5
$ws.Range("K135").Value = 12727.0386
$ws.Range("L135").Value = 13036.5
$ws.Range("M135").Value = -10192.0386
$ws.Range("N135").Value = -18106.5
$ws.Range("H137").Value = 12021.435
$ws.Range("I137").Value = 7348.3706
$ws.Range("J137").Value = 18662.105
$ws.Range("K137").Value = 22045.1118
$ws.Range("L137").Value = 55986.315
$ws.Range("M137").Value = -19495.1118
$ws.Range("N137").Value = -61086.315

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3026.7322
$ws.Range("I2").Value = 2202.319
$ws.Range("K2").Value = 2202.319
$ws.Range("M2").Value = -2089.319
$ws.Range("H4").Value = 399.5
$ws.Range("I4").Value = 359.2
$ws.Range("J4").Value = 439.8
$ws.Range("K4").Value = 359.2
$ws.Range("L4").Value = 439.8
$ws.Range("M4").Value = -243.2
$ws.Range("N4").Value = -671.8
$ws.Range("H5").Value = 145.66667
$ws.Range("I5").Value = 145.66667
$ws.Range("K5").Value = 145.66667
$ws.Range("M5").Value = -33.66667000000001
$ws.Range("H74").Value = 2657.88
$ws.Range("I74").Value = 995.62067
$ws.Range("J74").Value = 4953.381
$ws.Range("K74").Value = 995.62067
$ws.Range("L74").Value = 4953.381
$ws.Range("M74").Value = -121.62067
$ws.Range("N74").Value = -6701.381
$ws.Range("H77").Value = 2657.88
$ws.Range("I77").Value = 995.62067
$ws.Range("J77").Value = 4953.381
$ws.Range("K77").Value = 4978.10335
$ws.Range("L77").Value = 24766.905
$ws.Range("M77").Value = -610.1033500000003
$ws.Range("N77").Value = -33502.905
$ws.Range("H116").Value = 3026.7322
$ws.Range("I116").Value = 2202.319
$ws.Range("K116").Value = 2202.319
$ws.Range("M116").Value = 91.68100000000004

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3026.7322
$ws.Range("I3").Value = 2202.319
$ws.Range("K3").Value = 2202.319
$ws.Range("M3").Value = -2088.319
$ws.Range("H4").Value = 145.66667
$ws.Range("I4").Value = 145.66667
$ws.Range("K4").Value = 145.66667
$ws.Range("M4").Value = -30.66667000000001
$ws.Range("H26").Value = 9362.75
$ws.Range("I26").Value = 9362.75
$ws.Range("K26").Value = 9362.75
$ws.Range("M26").Value = -9070.75
$ws.Range("H96").Value = 15773
$ws.Range("I96").Value = 15773
$ws.Range("K96").Value = 15773
$ws.Range("M96").Value = -13027
$ws.Range("H134").Value = 8612.429
$ws.Range("I134").Value = 6234.4287
$ws.Range("J134").Value = 18124.428
$ws.Range("K134").Value = 18703.2861
$ws.Range("L134").Value = 54373.284
$ws.Range("M134").Value = -16168.2861
$ws.Range("N134").Value = -59443.284

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 31888.7
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 31888.7
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 31888.7
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -33138.7
$ws.Range("H51").Value = 28246.334
$ws.Range("I51").Value = 14745
$ws.Range("K51").Value = 14745
$ws.Range("M51").Value = -14009
$ws.Range("H58").Value = 4649.143
$ws.Range("I58").Value = 2046.4706
$ws.Range("K58").Value = 2046.4706
$ws.Range("M58").Value = -1843.4706
$ws.Range("H61").Value = 28246.334
$ws.Range("I61").Value = 14745
$ws.Range("K61").Value = 14745
$ws.Range("M61").Value = -14397
$ws.Range("H68").Value = 63565.43
$ws.Range("J68").Value = 63565.43
$ws.Range("L68").Value = 63565.43
$ws.Range("N68").Value = -65063.43
$ws.Range("H71").Value = 63565.43
$ws.Range("J71").Value = 63565.43
$ws.Range("L71").Value = 190696.29
$ws.Range("N71").Value = -198184.29
$ws.Range("H136").Value = 4649.143
$ws.Range("I136").Value = 2046.4706
$ws.Range("K136").Value = 6139.4118
$ws.Range("M136").Value = -3589.4118

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1969.9811
$ws.Range("I5").Value = 1072.1562
$ws.Range("J5").Value = 3338.0952
$ws.Range("K5").Value = 3216.4686
$ws.Range("L5").Value = 10014.2856
$ws.Range("M5").Value = -3104.4686
$ws.Range("N5").Value = -10238.2856
$ws.Range("H37").Value = 99953
$ws.Range("J37").Value = 99953
$ws.Range("L37").Value = 299859
$ws.Range("N37").Value = -300083
$ws.Range("H68").Value = 3166.9
$ws.Range("J68").Value = 3296.5557
$ws.Range("L68").Value = 9889.667099999999
$ws.Range("N68").Value = -11511.6671
$ws.Range("H71").Value = 3166.9
$ws.Range("J71").Value = 3296.5557
$ws.Range("L71").Value = 29669.0013
$ws.Range("N71").Value = -37781.0013
$ws.Range("H80").Value = 24159.727
$ws.Range("J80").Value = 15981.286
$ws.Range("L80").Value = 47943.858
$ws.Range("N80").Value = -49815.858
$ws.Range("H83").Value = 24159.727
$ws.Range("J83").Value = 15981.286
$ws.Range("L83").Value = 143831.574
$ws.Range("N83").Value = -153191.574
$ws.Range("H135").Value = 1969.9811
$ws.Range("I135").Value = 1072.1562
$ws.Range("J135").Value = 3338.0952
$ws.Range("K135").Value = 9649.405799999999
$ws.Range("L135").Value = 30042.8568
$ws.Range("M135").Value = -7114.405799999999
$ws.Range("N135").Value = -35112.8568

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 80000
$ws.Range("J103").Value = 80000
$ws.Range("L103").Value = 80000
$ws.Range("N103").Value = -82344

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 149680
$ws.Range("J116").Value = 149680
$ws.Range("L116").Value = 149680
$ws.Range("N116").Value = -158858

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 39987.5
$ws.Range("J68").Value = 39987.5
$ws.Range("L68").Value = 39987.5
$ws.Range("N68").Value = -41609.5
$ws.Range("H71").Value = 39987.5
$ws.Range("J71").Value = 39987.5
$ws.Range("L71").Value = 119962.5
$ws.Range("N71").Value = -128074.5
$ws.Range("H132").Value = 14381.857
$ws.Range("I132").Value = 11850.108
$ws.Range("J132").Value = 22188.084
$ws.Range("K132").Value = 35550.324
$ws.Range("L132").Value = 66564.25199999999
$ws.Range("M132").Value = -33020.324
$ws.Range("N132").Value = -71624.25199999999

Write-Output "Applied all Twintania_Profits updates"